$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the format of G1 (bold, centered, bordered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the data values for the new column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
